# Atualização de bases das ligas, do dia: 21-04-2024 às 14:32
#
# This script applies two logical changes to the single worksheet
# ("Mexico Liga de Expansion"):
#   1. Rows 91 and 92 had their match data swapped (the row-index
#      column A stays put, everything else - id, odds, results - moves).
#   2. A brand new match result (row 234) is appended at the bottom of
#      the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Swap the contents of row 91 and row 92 (columns B..AC). Column A
#    (the sequential id column) is intentionally left untouched.
# ---------------------------------------------------------------------

# New content for row 91 (this used to be row 92's data)
$ws.Range("B91").Value = 6924569
$ws.Range("C91").Value = "Mexico Liga de Expansion"
$ws.Range("D91").Value = "Mexico Liga de Expansion"
$ws.Range("F91").Value = "Venados FC"
$ws.Range("G91").Value = "Dorados"
$ws.Range("H91").Value = 4
$ws.Range("I91").Value = 1
$ws.Range("J91").Value = "H"
$ws.Range("K91").Value = 1.615
$ws.Range("L91").Value = 4
$ws.Range("M91").Value = 4.5
$ws.Range("N91").Value = 1.5
$ws.Range("O91").Value = 4.75
$ws.Range("P91").Value = 5.75
$ws.Range("Q91").Value = -1.25
$ws.Range("R91").Value = 1.925
$ws.Range("S91").Value = 1.875
$ws.Range("T91").Value = 3
$ws.Range("U91").Value = 1.75
$ws.Range("V91").Value = 1.95
$ws.Range("W91").Value = 0.5
$ws.Range("X91").Value = -1
$ws.Range("Y91").Value = -1
$ws.Range("Z91").Value = 0.925
$ws.Range("AA91").Value = -1
$ws.Range("AB91").Value = 0.75
$ws.Range("AC91").Value = -1

# New content for row 92 (this used to be row 91's data)
$ws.Range("B92").Value = 6924568
$ws.Range("C92").Value = "Mexico Liga de Expansion"
$ws.Range("D92").Value = "Mexico Liga de Expansion"
$ws.Range("F92").Value = "Atletico Morelia"
$ws.Range("G92").Value = "Atlante"
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 1
$ws.Range("J92").Value = "A"
$ws.Range("K92").Value = 2.4
$ws.Range("L92").Value = 3
$ws.Range("M92").Value = 2.875
$ws.Range("N92").Value = 2.7
$ws.Range("O92").Value = 3.1
$ws.Range("P92").Value = 2.8
$ws.Range("Q92").Value = 0
$ws.Range("R92").Value = 1.85
$ws.Range("S92").Value = 1.95
$ws.Range("T92").Value = 2.25
$ws.Range("U92").Value = 1.975
$ws.Range("V92").Value = 1.725
$ws.Range("W92").Value = -1
$ws.Range("X92").Value = -1
$ws.Range("Y92").Value = 1.8
$ws.Range("Z92").Value = -1
$ws.Range("AA92").Value = 0.95
$ws.Range("AB92").Value = -1
$ws.Range("AC92").Value = 0.7250000000000001

# ---------------------------------------------------------------------
# 2) Append the new row 234 at the bottom of the sheet.
#    Column A uses the same bold/boxed style as the other id cells
#    (copied from A91) and column E uses the same date/time number
#    format as the other date cells (copied from E91).
# ---------------------------------------------------------------------

$ws.Range("A91").Copy()
$ws.Range("A234").PasteSpecial(-4122)
$ws.Range("E91").Copy()
$ws.Range("E234").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A234").Value = 232
$ws.Range("B234").Value = 8117254
$ws.Range("C234").Value = "Mexico Liga de Expansion"
$ws.Range("D234").Value = "Mexico Liga de Expansion"
$ws.Range("E234").Value = 45403
$ws.Range("F234").Value = "Club Atletico La Paz"
$ws.Range("G234").Value = "Oaxaca"
$ws.Range("H234").Value = 3
$ws.Range("I234").Value = 2
$ws.Range("J234").Value = "H"
$ws.Range("K234").Value = 1.909
$ws.Range("L234").Value = 3.5
$ws.Range("M234").Value = 3.4
$ws.Range("N234").Value = 1.85
$ws.Range("O234").Value = 3.5
$ws.Range("P234").Value = 3.5
$ws.Range("Q234").Value = -0.5
$ws.Range("R234").Value = 1.85
$ws.Range("S234").Value = 1.95
$ws.Range("T234").Value = 2.5
$ws.Range("U234").Value = 1.775
$ws.Range("V234").Value = 1.925
$ws.Range("W234").Value = 0.8500000000000001
$ws.Range("X234").Value = -1
$ws.Range("Y234").Value = -1
$ws.Range("Z234").Value = 0.8500000000000001
$ws.Range("AA234").Value = -1
$ws.Range("AB234").Value = 0.7749999999999999
$ws.Range("AC234").Value = -1
